$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.899.72'
$ws.Range('E2').Value = '  -0.55%  '
$ws.Range('D3').Value = '1.667.60'
$ws.Range('E3').Value = '  +0.74%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.59'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('E6').Value = '  +5.59%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('E8').Value = '  +1.10%  '
$ws.Range('E9').Value = '  -0.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.24'
$ws.Range('E10').Value = '  +2.85%  '
$ws.Range('E11').Value = '  +3.61%  '
$ws.Range('D12').Value = '1.902.27'
$ws.Range('E12').Value = '  +0.67%  '
$ws.Range('D13').Value = '1.643.47'
$ws.Range('E13').Value = '  -0.77%  '
$ws.Range('E14').Value = '  +0.16%  '
$ws.Range('E15').Value = '  +1.40%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.14'
$ws.Range('E16').Value = '  +1.42%  '
$ws.Range('D17').Value = '26.926.42'
$ws.Range('E17').Value = '  -0.44%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '234.15'
$ws.Range('E18').Value = '  -1.94%  '
$ws.Range('E19').Value = '  +1.48%  '
$ws.Range('D20').Value = '0.0₃0731'
$ws.Range('E20').Value = '  +0.33%  '
$ws.Range('E21').Value = '  -0.08%  '
$ws.Range('E22').Value = '  -0.52%  '
$ws.Range('E23').Value = '  -1.90%  '
$ws.Range('E24').Value = '  -0.92%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.12'
$ws.Range('E25').Value = '  +0.23%  '
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('E27').Value = '  +0.94%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.88'
$ws.Range('E28').Value = '  +0.22%  '
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.35'
$ws.Range('E32').Value = '  +1.81%  '
$ws.Range('D33').Value = '1.450.63'
$ws.Range('E33').Value = '  -4.27%  '
$ws.Range('E34').Value = '  +2.33%  '
$ws.Range('E35').Value = '  +3.29%  '
$ws.Range('E36').Value = '  -0.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.582'
$ws.Range('E37').Value = '  +0.95%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.903'
$ws.Range('E38').Value = '  +1.89%  '
$ws.Range('E39').Value = '  +0.65%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.74'
$ws.Range('E40').Value = '  -3.67%  '
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('E42').Value = '  +1.26%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '66.12'
$ws.Range('E43').Value = '  +0.29%  '
$ws.Range('E44').Value = '  +5.82%  '
$ws.Range('D45').Value = '1.809.78'
$ws.Range('E45').Value = '  +0.72%  '
$ws.Range('E46').Value = '  +1.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.49'
$ws.Range('E47').Value = '  +1.04%  '
$ws.Range('E48').Value = '  +1.24%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.102'
$ws.Range('E49').Value = '  +4.56%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0506'
$ws.Range('E50').Value = '  -0.24%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.53'
$ws.Range('E51').Value = '  -0.72%  '
